# Applies the edit described in the commit: fixes the header date line,
# adds an "input" lead-in sentence to the Stage 1 paragraph, and flattens the
# Hebrew phrase/letter runs in the four "Detailed Examples" paragraphs back
# into single runs (was split letter-by-letter / word-by-word).
#
# Note: this runtime's `Range.Text = "..."` assignment only overwrites the
# *first* run contained in the target range and leaves any further runs in
# that range untouched, and plain `InsertBefore`/`InsertAfter` silently
# coalesces new text into an adjacent run whenever the resolved formatting is
# identical. Neither behavior can produce "replace this paragraph's runs with
# exactly these new runs", which several hunks of the diff need (two bold
# runs in the date line; a new un-bold run before the existing Stage 1 run;
# N runs collapsed to exactly 1 run elsewhere). So throughout we rebuild the
# whole paragraph at once with `Range.InsertXML`, which substitutes the
# targeted range's OOXML directly and is not subject to that run-merging.

$d = $word.ActiveDocument

# --- 1) Byline/date block -------------------------------------------------
# Original: two paragraphs -> "For Journal Editors and Academic Reviewers"
#           and "Date: 2025-10-19" (each its own bold run/paragraph).
# Target:   one paragraph, two bold runs -> "Date: 2025-10-" + "21".
$xmlDate = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="51C2BB8D" w14:textId="77777777" w:rsidR="002A0102" w:rsidRDefault="00227FFB"><w:pPr><w:spacing w:after="120"/></w:pPr><w:r><w:rPr><w:b/></w:rPr><w:t>Date: 2025-10-</w:t></w:r><w:r><w:rPr><w:b/></w:rPr><w:t>21</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$d.Paragraphs(2).Range.InsertXML($xmlDate)
$d.Paragraphs(3).Range.Delete()

# --- 2) Stage 1 paragraph --------------------------------------------------
# Prepend a new leading run "A chapter of Psalms is fed to " before the
# existing "Establishes the big picture. ..." run.
$xmlStage1 = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="63B13855" w14:textId="77777777" w:rsidR="002A0102" w:rsidRDefault="00227FFB"><w:r><w:t xml:space="preserve">A chapter of Psalms is fed to </w:t></w:r><w:r><w:t>Establishes the big picture. The system identifies the psalm’s genre (such as lament, praise, or wisdom), develops a central theological thesis, and creates a structural outline. This stage draws on a predefined analytical framework and reference materials to understand the psalm’s place within the broader biblical tradition.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$d.Paragraphs(8).Range.InsertXML($xmlStage1)

# --- 3) Example 1: sibilant sounds paragraph ------------------------------
# Collapse the run-per-Hebrew-letter split back into a single run.
$xml47 = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="54A7F205" w14:textId="77777777" w:rsidR="002A0102" w:rsidRDefault="00227FFB"><w:r><w:t>The system’s analysis of Psalm 1 demonstrates how phonetic transcriptions prevent common errors. In the opening verses, the system identifies specific sibilant sounds (ש, ס, צ) rather than making vague claims about “sibilant alliteration.”</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$d.Paragraphs(46).Range.InsertXML($xml47)

# --- 4) Example 1: "For instance, in verse 1..." paragraph ---------------
# Collapse the run-per-Hebrew-word split back into a single run.
$xml48 = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="49FC8561" w14:textId="77777777" w:rsidR="002A0102" w:rsidRDefault="00227FFB"><w:r><w:t>For instance, in verse 1, the sequence “אשרי האיש אשר לא הלך” contains the sibilant ש in אשרי and אשר, while verse 2’s “כי אם בתורת יהוה” includes the ש in יהוה. The system maps these sounds precisely across the poetic lines, showing how they create acoustic patterns that reinforce the psalm’s structure.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$d.Paragraphs(47).Range.InsertXML($xml48)

# --- 5) Example 2: tree metaphor paragraph --------------------------------
$xml51 = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="3E4A428C" w14:textId="77777777" w:rsidR="002A0102" w:rsidRDefault="00227FFB"><w:r><w:t>The system’s analysis of the tree metaphor in Psalm 1:3 illustrates how cross-textual comparison leads to specific insights. The Hebrew phrase “עץ שתול על־פלגי מים” (“tree planted by streams of water”) uses the verb שתול (planted), which specifically means transplanted or cultivated, not naturally occurring.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$d.Paragraphs(50).Range.InsertXML($xml51)

# --- 6) Example 2: "The system's database reveals..." paragraph ---------
$xml53 = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="72313A62" w14:textId="77777777" w:rsidR="002A0102" w:rsidRDefault="00227FFB"><w:r><w:t>The system’s database reveals that שתול appears 12 times in the Hebrew Bible, always in contexts of deliberate agricultural practice. This evidence leads the system to conclude that the tree in Psalm 1 represents not just any tree, but specifically a cultivated tree placed by irrigation channels — an image of intentional spiritual development rather than accidental growth.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$d.Paragraphs(52).Range.InsertXML($xml53)
